# Apply updated Leve profit-calculation figures (currentAveragePrice / LevePrice / LeveProfit
# columns H:N) across the ALC, ARM, BSM, CRP, CUL and WVR sheets, per the latest market-board
# pull from the scheduled runner.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 62: The Mustache Suits Him
$ws.Range("H62").Value = 1531.3334
$ws.Range("I62").Value = 1531.3334
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 1531.3334
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -907.3334
$ws.Range("N62").Value = ""

# Row 65: Forgery of Convenience (L)
$ws.Range("H65").Value = 1531.3334
$ws.Range("I65").Value = 1531.3334
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 7656.666999999999
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -4536.666999999999
$ws.Range("N65").Value = ""

# Row 86: Filling in the Blanks
$ws.Range("H86").Value = 2275.5293
$ws.Range("I86").Value = 1754.2858
$ws.Range("J86").Value = 2640.4
$ws.Range("K86").Value = 1754.2858
$ws.Range("L86").Value = 2640.4
$ws.Range("M86").Value = -631.2858000000001
$ws.Range("N86").Value = -4886.4

# Row 89: Ink into Antiquity (L)
$ws.Range("H89").Value = 2275.5293
$ws.Range("I89").Value = 1754.2858
$ws.Range("J89").Value = 2640.4
$ws.Range("K89").Value = 8771.429
$ws.Range("L89").Value = 13202
$ws.Range("M89").Value = -3155.429
$ws.Range("N89").Value = -24434

# Row 125: Body over Mind
$ws.Range("H125").Value = 7455.5
$ws.Range("I125").Value = 10000
$ws.Range("J125").Value = 6365
$ws.Range("K125").Value = 90000
$ws.Range("L125").Value = 57285
$ws.Range("M125").Value = -87540
$ws.Range("N125").Value = -62205

# Row 132: Fast-forwarding Flora
$ws.Range("H132").Value = 1526157
$ws.Range("I132").Value = 1690062
$ws.Range("K132").Value = 5070186
$ws.Range("M132").Value = -5067656

$ws = $wb.Worksheets.Item("ARM")
# Row 32: Ingot We Trust
$ws.Range("H32").Value = 1093.43
$ws.Range("I32").Value = 948.5484
$ws.Range("J32").Value = 3018.2856
$ws.Range("K32").Value = 948.5484
$ws.Range("L32").Value = 3018.2856
$ws.Range("M32").Value = -661.5484
$ws.Range("N32").Value = -3592.2856

# Row 74: As the Bolt Flies
$ws.Range("H74").Value = 749.6539
$ws.Range("I74").Value = 619.3684
$ws.Range("J74").Value = 1103.2858
$ws.Range("K74").Value = 619.3684
$ws.Range("L74").Value = 1103.2858
$ws.Range("M74").Value = 254.6316
$ws.Range("N74").Value = -2851.2858

# Row 77: Heavy Metal Banned (L)
$ws.Range("H77").Value = 749.6539
$ws.Range("I77").Value = 619.3684
$ws.Range("J77").Value = 1103.2858
$ws.Range("K77").Value = 3096.842
$ws.Range("L77").Value = 5516.429
$ws.Range("M77").Value = 1271.158
$ws.Range("N77").Value = -14252.429

# Row 132: Don't Bore Me, Ore Me
$ws.Range("H132").Value = 1361.3715
$ws.Range("I132").Value = 1107.6129
$ws.Range("J132").Value = 3328
$ws.Range("K132").Value = 3322.8387
$ws.Range("L132").Value = 9984
$ws.Range("M132").Value = -792.8387000000002
$ws.Range("N132").Value = -15044

$ws = $wb.Worksheets.Item("BSM")
# Row 80: Unbreaker
$ws.Range("H80").Value = 838.58826
$ws.Range("I80").Value = 260
$ws.Range("J80").Value = 962.5714
$ws.Range("K80").Value = 260
$ws.Range("L80").Value = 962.5714
$ws.Range("M80").Value = 738
$ws.Range("N80").Value = -2958.5714

# Row 83: Attack on Titanium (L)
$ws.Range("H83").Value = 838.58826
$ws.Range("I83").Value = 260
$ws.Range("J83").Value = 962.5714
$ws.Range("K83").Value = 1300
$ws.Range("L83").Value = 4812.857
$ws.Range("M83").Value = 3692
$ws.Range("N83").Value = -14796.857

# Row 99: Meddle in Metal
$ws.Range("H99").Value = 1410.28
$ws.Range("I99").Value = 1337.3158
$ws.Range("J99").Value = 1641.3334
$ws.Range("K99").Value = 1337.3158
$ws.Range("L99").Value = 1641.3334
$ws.Range("M99").Value = 160.6841999999999
$ws.Range("N99").Value = -4637.3334

# Row 107: The Gold Experience
$ws.Range("H107").Value = 1317.3846
$ws.Range("I107").Value = 1330.08
$ws.Range("J107").Value = 1000
$ws.Range("K107").Value = 1330.08
$ws.Range("L107").Value = 1000
$ws.Range("M107").Value = 589.9200000000001
$ws.Range("N107").Value = -4840

# Row 134: Ruthenium Supremium
$ws.Range("H134").Value = 1791.6
$ws.Range("I134").Value = 1229.25
$ws.Range("J134").Value = 2791.3333
$ws.Range("K134").Value = 3687.75
$ws.Range("L134").Value = 8373.999899999999
$ws.Range("M134").Value = -1152.75
$ws.Range("N134").Value = -13443.9999

$ws = $wb.Worksheets.Item("CRP")
# Row 31: Wall Not Found
$ws.Range("H31").Value = 1691.4117
$ws.Range("I31").Value = 1400.9524
$ws.Range("J31").Value = 3046.889
$ws.Range("K31").Value = 1400.9524
$ws.Range("L31").Value = 3046.889
$ws.Range("M31").Value = -1105.9524
$ws.Range("N31").Value = -3636.889

# Row 34: Armoires of the Rich and Famous
$ws.Range("H34").Value = 1691.4117
$ws.Range("I34").Value = 1400.9524
$ws.Range("J34").Value = 3046.889
$ws.Range("K34").Value = 1400.9524
$ws.Range("L34").Value = 3046.889
$ws.Range("M34").Value = -1198.9524
$ws.Range("N34").Value = -3450.889

# Row 96: Composition
$ws.Range("H96").Value = 24606.715
$ws.Range("J96").Value = 24606.715
$ws.Range("L96").Value = 24606.715
$ws.Range("N96").Value = -30098.715

# Row 99: O Pine
$ws.Range("H99").Value = 2426.7
$ws.Range("J99").Value = 2155.8572
$ws.Range("L99").Value = 2155.8572
$ws.Range("N99").Value = -5151.8572

# Row 119: Off to a Good Staff
$ws.Range("H119").Value = 39530.5
$ws.Range("J119").Value = 39530.5
$ws.Range("L119").Value = 39530.5
$ws.Range("N119").Value = -49206.5

# Row 126: A Better Conductor
$ws.Range("H126").Value = 2426.7
$ws.Range("J126").Value = 2155.8572
$ws.Range("L126").Value = 6467.571599999999
$ws.Range("N126").Value = -11407.5716

# Row 132: Hull Lotta Damage
$ws.Range("H132").Value = 2801.111
$ws.Range("I132").Value = 2102
$ws.Range("K132").Value = 6306
$ws.Range("M132").Value = -3776

# Row 134: Wood You Be Quiet
$ws.Range("H134").Value = 1131.1714
$ws.Range("I134").Value = 963.12
$ws.Range("J134").Value = 1551.3
$ws.Range("K134").Value = 2889.36
$ws.Range("L134").Value = 4653.9
$ws.Range("M134").Value = -354.3600000000001
$ws.Range("N134").Value = -9723.9

$ws = $wb.Worksheets.Item("CUL")
# Row 55: Pagan Pastries
$ws.Range("H55").Value = 26405.562
$ws.Range("I55").Value = 156175.6
$ws.Range("J55").Value = 2374.074
$ws.Range("K55").Value = 468526.8
$ws.Range("L55").Value = 7122.222
$ws.Range("M55").Value = -468349.8
$ws.Range("N55").Value = -7476.222

# Row 109: Cure for What Ails
$ws.Range("H109").Value = 3358.75
$ws.Range("I109").Value = 133.33333
$ws.Range("J109").Value = 5294
$ws.Range("K109").Value = 399.99999
$ws.Range("L109").Value = 15882
$ws.Range("M109").Value = 640.00001
$ws.Range("N109").Value = -17962

# Row 121: A Cookie for Your Troubles
$ws.Range("H121").Value = 4763.615
$ws.Range("I121").Value = 250
$ws.Range("J121").Value = 5139.75
$ws.Range("K121").Value = 750
$ws.Range("L121").Value = 15419.25
$ws.Range("M121").Value = 560
$ws.Range("N121").Value = -18039.25

$ws = $wb.Worksheets.Item("WVR")
# Row 132: Comfy Cabins
$ws.Range("H132").Value = 1644.762
$ws.Range("I132").Value = 1304.1333
$ws.Range("J132").Value = 2496.3333
$ws.Range("K132").Value = 3912.3999
$ws.Range("L132").Value = 7488.999899999999
$ws.Range("M132").Value = -1382.3999
$ws.Range("N132").Value = -12548.9999

# Row 136: Weaving the Envelope
$ws.Range("H136").Value = 4295.9375
$ws.Range("I136").Value = 944.36365
$ws.Range("J136").Value = 11669.4
$ws.Range("K136").Value = 2833.09095
$ws.Range("L136").Value = 35008.2
$ws.Range("M136").Value = -283.0909499999998
$ws.Range("N136").Value = -40108.2
